# "Generate Report for Handback"
# Marks both localized files (zh-cn, de-de) as handed back: populates the
# "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime"
# columns on each language sheet, and flips the Overview status text from
# "Ready for handoff" to "Handed back: in sync with en-US".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Overview sheet: status text for both files/languages
# ---------------------------------------------------------------------
$ovw = $wb.Worksheets.Item("Overview")

$ovw.Range("E2").Value = "Handed back: in sync with en-US"
$ovw.Range("F2").Value = "Handed back: in sync with en-US"
$ovw.Range("E3").Value = "Handed back: in sync with en-US"
$ovw.Range("F3").Value = "Handed back: in sync with en-US"

$ovw.Columns.Item(5).ColumnWidth = 29.1666666666667
$ovw.Columns.Item(6).ColumnWidth = 29.1666666666667

# ---------------------------------------------------------------------
# 2) zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Columns.Item(3).ColumnWidth = 29.1666666666667
$zh.Columns.Item(9).ColumnWidth = 39.1666666666667
$zh.Columns.Item(10).ColumnWidth = 39.1666666666667

$zh.Hyperlinks.Add(
    $zh.Range("I2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e004a855ce08326f1f5a5644f127a3f5f8441729/e2e/77879806-9de0-495c-b6c6-b7169e833960.md",
    "",
    "",
    "77879806-9de0-495c-b6c6-b7169e833960.md"
) | Out-Null
$zh.Range("J2").Value = "77879806-9de0-495c-b6c6-b7169e833960.bee9ce64a15931c20b8ced5c1de9ba1434e61f1a.zh-cn.xlf"
$zh.Range("K2").Value = "2016-09-06 23:16:31"

$zh.Hyperlinks.Add(
    $zh.Range("I3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e004a855ce08326f1f5a5644f127a3f5f8441729/e2e/87db403b-203f-44e4-a0e0-0276bc3326ca.md",
    "",
    "",
    "87db403b-203f-44e4-a0e0-0276bc3326ca.md"
) | Out-Null
$zh.Range("J3").Value = "87db403b-203f-44e4-a0e0-0276bc3326ca.e648f48d8489fd8853d6d56454fe0e325d8811ee.zh-cn.xlf"
$zh.Range("K3").Value = "2016-09-06 23:16:31"

# ---------------------------------------------------------------------
# 3) de-de sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Columns.Item(3).ColumnWidth = 29.1666666666667
$de.Columns.Item(9).ColumnWidth = 39.1666666666667
$de.Columns.Item(10).ColumnWidth = 39.1666666666667

$de.Hyperlinks.Add(
    $de.Range("I2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e004a855ce08326f1f5a5644f127a3f5f8441729/e2e/77879806-9de0-495c-b6c6-b7169e833960.md",
    "",
    "",
    "77879806-9de0-495c-b6c6-b7169e833960.md"
) | Out-Null
$de.Range("J2").Value = "77879806-9de0-495c-b6c6-b7169e833960.bee9ce64a15931c20b8ced5c1de9ba1434e61f1a.de-de.xlf"
$de.Range("K2").Value = "2016-09-06 23:16:40"

$de.Hyperlinks.Add(
    $de.Range("I3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e004a855ce08326f1f5a5644f127a3f5f8441729/e2e/87db403b-203f-44e4-a0e0-0276bc3326ca.md",
    "",
    "",
    "87db403b-203f-44e4-a0e0-0276bc3326ca.md"
) | Out-Null
$de.Range("J3").Value = "87db403b-203f-44e4-a0e0-0276bc3326ca.e648f48d8489fd8853d6d56454fe0e325d8811ee.de-de.xlf"
$de.Range("K3").Value = "2016-09-06 23:16:40"
